$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: drop the trailing empty placeholder cells (I20, K20:R20) ---
$ws.Range("I20").ClearContents()
$ws.Range("K20:R20").ClearContents()

# --- Row 21: new product line ---
$ws.Range("A21").Value = "6XS18316"
$ws.Range("B21").Value = "B.GREEN OMEGA 3 1000MG 30 CAPS"
$ws.Range("D21").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E21").Value = "Tiene ES"
$ws.Range("F21").Value = "No Tiene IT - TRADOTTO"

# G21 ("Cantidad Neta") is stored as text in this sheet (like every other
# row), so force text formatting, write the value, then drop the format
# again so the cell keeps the default (unstyled) look of its neighbours.
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "30"
$ws.Range("G21").ClearFormats()

$ws.Range("H21").Value = "UND"
$ws.Range("J21").Value = "Revisado y Traducido"
